$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.625.60'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  +1.72%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''3.465.06'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  -0.04%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = '''  -0.18%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''581.73'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  +1.27%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''167.99'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  +4.64%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').Value = '''  -0.14%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = '''3.462.99'
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Value = '''0.567'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '''  -1.10%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = '''7.28'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  +0.80%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('E11').Value = '''  +1.35%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('E12').Value = '''  -1.54%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = '''4.056.22'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '''  -0.20%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('E14').Value = '''  +0.06%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''27.49'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  -0.42%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('E16').Value = '''  +0.16%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''65.527.63'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  +1.27%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''3.468.75'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  +0.18%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('E19').Value = '''  +0.00%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''13.77'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  -0.73%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = '''385.48'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '''  +1.07%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''7.93'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  -0.27%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('E23').Value = '''  +0.11%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = '''71.52'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '''  -1.51%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''0.521'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  -1.58%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('E26').Value = '''  +1.16%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = '''9.80'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '''  -0.51%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = '''0.181'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '''  +1.40%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = '''0.995'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '''  -0.59%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = '''6.24'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '''  +1.89%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('E31').Value = '''  +0.50%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''2.03'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  +0.87%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = '''23.31'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '''  -0.09%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''7.32'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  +3.92%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('E35').Value = '''  +0.04%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = '''1.52'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '''  -3.80%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = '''160.19'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '''  -0.44%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = '''0.892'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '''  +8.26%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('E39').Value = '''  -0.59%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = '''0.0735'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  -1.49%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('B41').Value = '''RenderToken'
$ws.Range('B41').Style = "Normal"
$ws.Range('C41').Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C41').Style = "Normal"
$ws.Range('D41').Value = '''6.62'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '''  +2.41%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('B42').Value = '''EnergySwap'
$ws.Range('B42').Style = "Normal"
$ws.Range('C42').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C42').Style = "Normal"
$ws.Range('D42').Value = '''26.09'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''  -2.94%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('B43').Value = '''InjectiveProtocol'
$ws.Range('B43').Style = "Normal"
$ws.Range('C43').Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C43').Style = "Normal"
$ws.Range('D43').Value = '''26.82'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '''  +3.69%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('B44').Value = '''Maker'
$ws.Range('B44').Style = "Normal"
$ws.Range('C44').Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C44').Style = "Normal"
$ws.Range('D44').Value = '''2.798.88'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''  -1.37%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = '''43.06'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''  +0.60%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = '''4.46'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '''  -1.05%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = '''0.0309'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '''  -0.40%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = '''  +2.82%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = '''337.04'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '''  +0.69%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = '''1.07'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '''  +1.39%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = '''32.46'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '''  +5.13%  '
$ws.Range('E51').Style = "Normal"
